$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(234, 44308, 0, 12, 145.8434613514827),
    @(235, 44309, 1, 10, 121.5362177929023),
    @(236, 44310, 0, 8, 97.22897423432183),
    @(237, 44311, 3, 8, 97.22897423432183),
    @(238, 44312, 6, 13, 157.997083130773)
)

foreach ($row in $data) {
    $r = $row[0]
    $prev = $r - 1

    # Copy formatting (style) of column A from the previous row's date cell
    $ws.Cells.Item($prev, 1).Copy($ws.Cells.Item($r, 1))

    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}
